$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.769.05'
$ws.Range("E2").Value = '  -0.42%  '
$ws.Range("D3").Value = '2.252.93'
$ws.Range("E3").Value = '  -1.04%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '''303.70'
$ws.Range("E5").Value = '  -0.65%  '
$ws.Range("D6").Value = '''94.61'
$ws.Range("E6").Value = '  +1.50%  '
$ws.Range("D7").Value = '''0.523'
$ws.Range("E7").Value = '  -1.37%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("E9").Value = '  -0.56%  '
$ws.Range("D10").Value = '''34.55'
$ws.Range("E10").Value = '  +5.28%  '
$ws.Range("E11").Value = '  -1.66%  '
$ws.Range("E12").Value = '  -0.42%  '
$ws.Range("E13").Value = '  -1.66%  '
$ws.Range("D14").Value = '2.608.58'
$ws.Range("E14").Value = '  -0.74%  '
$ws.Range("D15").Value = '''14.25'
$ws.Range("E15").Value = '  -0.94%  '
$ws.Range("D16").Value = '2.258.27'
$ws.Range("E16").Value = '  -0.65%  '
$ws.Range("D17").Value = '''0.786'
$ws.Range("E17").Value = '  +0.15%  '
$ws.Range("D18").Value = '41.671.10'
$ws.Range("E18").Value = '  -0.49%  '
$ws.Range("E19").Value = '  -4.86%  '
$ws.Range("D20").Value = '0.0₃0895'
$ws.Range("E20").Value = '  -2.33%  '
$ws.Range("E21").Value = '  -0.99%  '
$ws.Range("D22").Value = '''67.77'
$ws.Range("E22").Value = '  -0.39%  '
$ws.Range("D23").Value = '''236.42'
$ws.Range("E23").Value = '  -3.06%  '
$ws.Range("E24").Value = '  -2.32%  '
$ws.Range("D25").Value = '''0.999'
$ws.Range("E25").Value = '  -0.23%  '
$ws.Range("E26").Value = '  -1.61%  '
$ws.Range("E27").Value = '  -2.11%  '
$ws.Range("D28").Value = '''36.18'
$ws.Range("E28").Value = '  +3.05%  '
$ws.Range("E29").Value = '  +1.32%  '
$ws.Range("D30").Value = '''9.42'
$ws.Range("E30").Value = '  -2.97%  '
$ws.Range("D31").Value = '''159.28'
$ws.Range("E31").Value = '  +0.06%  '
$ws.Range("D32").Value = '''0.999'
$ws.Range("E32").Value = '  -0.01%  '
$ws.Range("E33").Value = '  -4.02%  '
$ws.Range("E34").Value = '  +2.72%  '
$ws.Range("D35").Value = '''0.0731'
$ws.Range("E35").Value = '  -1.85%  '
$ws.Range("D36").Value = '''16.98'
$ws.Range("E36").Value = '  -2.20%  '
$ws.Range("E37").Value = '  +0.54%  '
$ws.Range("E38").Value = '  -1.08%  '
$ws.Range("E39").Value = '  +0.55%  '
$ws.Range("E40").Value = '  -2.66%  '
$ws.Range("D41").Value = '''3.97'
$ws.Range("E41").Value = '  +0.50%  '
$ws.Range("D42").Value = '''2.36'
$ws.Range("E42").Value = '  +4.27%  '
$ws.Range("D43").Value = '1.970.91'
$ws.Range("E43").Value = '  -1.92%  '
$ws.Range("D44").Value = '''0.0281'
$ws.Range("E44").Value = '  -0.70%  '
$ws.Range("D45").Value = '''18.70'
$ws.Range("E45").Value = '  -6.18%  '
$ws.Range("D46").Value = '''2.90'
$ws.Range("E46").Value = '  -0.81%  '
$ws.Range("E47").Value = '  -4.49%  '
$ws.Range("D48").Value = '''53.01'
$ws.Range("E48").Value = '  -1.11%  '
$ws.Range("D49").Value = '''72.44'
$ws.Range("E49").Value = '  -0.29%  '
$ws.Range("D50").Value = '''1.49'
$ws.Range("E50").Value = '  -1.27%  '
$ws.Range("D51").Value = '''90.40'
$ws.Range("E51").Value = '  -1.73%  '
